# Update "想去人数" (wish-to-go count) figures across the workbook sheets.
# Mirrors a data refresh where several events' interest counts increased.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 2437
$ws1.Range("F12").Value = 877
$ws1.Range("F13").Value = 1084
$ws1.Range("F21").Value = 1092
$ws1.Range("F25").Value = 209
$ws1.Range("F29").Value = 3016
$ws1.Range("F30").Value = 460
$ws1.Range("F36").Value = 1568
$ws1.Range("F40").Value = 128
$ws1.Range("F44").Value = 118
$ws1.Range("F45").Value = 77

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 4411
$ws2.Range("F14").Value = 160

# --- Sheet: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2247
$ws3.Range("F3").Value = 708

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 2247
$ws4.Range("F3").Value  = 708
$ws4.Range("F9").Value  = 2437
$ws4.Range("F14").Value = 877
$ws4.Range("F15").Value = 1084
$ws4.Range("F24").Value = 1092
$ws4.Range("F29").Value = 209
$ws4.Range("F32").Value = 3016
$ws4.Range("F34").Value = 460
$ws4.Range("F38").Value = 1568
$ws4.Range("F42").Value = 128
$ws4.Range("F47").Value = 77

$wb.Save()
